$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49, shifting the existing rows 49:166 down to 50:167
$ws.Rows("49:49").Insert()

# Populate the newly inserted row 49 with the new daily price record.
# The "constant" columns mirror every other data row in this table.
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 45012
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100108
$ws.Range("H49").Value = "Tropicales y subtropicales"
$ws.Range("I49").Value = 100108002
$ws.Range("J49").Value = "Mango"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 248
$ws.Range("N49").Value = 7000
$ws.Range("O49").Value = 7000
$ws.Range("P49").Value = 7000
$ws.Range("Q49").Value = "$/bandeja 4 kilos"
$ws.Range("R49").Value = "Perú"
$ws.Range("S49").Value = 1750
$ws.Range("T49").Value = 4
